$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataObat")

$ws.Range("B2").Value = "kapsul"
$ws.Range("C2").Value = "keras"
$ws.Range("D2").Value = 20000
$ws.Range("E2").Value = 30
$ws.Range("H2").Value = "28/06/2045"
